# Atualização de bases das ligas, do dia: 15-04-2024 às 22:35
# Swap the data (columns B..AC) between pairs of rows, keeping the
# "id" column (A) fixed per row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowPairs = @(
    @(8, 9),
    @(10, 11),
    @(22, 23),
    @(43, 44),
    @(56, 57)
)

foreach ($pair in $rowPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $range1 = $ws.Range("B" + $r1 + ":AC" + $r1)
    $range2 = $ws.Range("B" + $r2 + ":AC" + $r2)

    $values1 = $range1.Value2
    $values2 = $range2.Value2

    $range1.Value2 = $values2
    $range2.Value2 = $values1
}
